# Atualização de bases das ligas, do dia: 18-02-2024 às 22:54
#
# Azerbaijan Premier League sheet1 update:
#  - rows 82/83 (B:AC) swap places (two matches on same date got reordered)
#  - rows 84/85 (B:AC) swap places (two matches on same date got reordered)
#  - rows 200/201 get their final match data (result, odds, etc.) filled in
#    and are moved earlier in the sheet (the two matches that used to sit
#    at 200/201 are pushed down to new rows 202/203, now themselves filled
#    in with final data as well)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Swap rows 82 and 83 (columns B through AC; A/C/D/E are identical
#    between the two rows already, so they don't need touching)
# ---------------------------------------------------------------------
$tmp82 = $ws.Range("B82:AC82").Value()
$tmp83 = $ws.Range("B83:AC83").Value()
$ws.Range("B82:AC82").Value = $tmp83
$ws.Range("B83:AC83").Value = $tmp82

# ---------------------------------------------------------------------
# 2) Swap rows 84 and 85 (columns B through AC)
# ---------------------------------------------------------------------
$tmp84 = $ws.Range("B84:AC84").Value()
$tmp85 = $ws.Range("B85:AC85").Value()
$ws.Range("B84:AC84").Value = $tmp85
$ws.Range("B85:AC85").Value = $tmp84

# ---------------------------------------------------------------------
# 3) Row 200 now carries the final result for match id 7011603
#    (previously this row held the not-yet-played 7011605 fixture)
# ---------------------------------------------------------------------
$ws.Range("B200").Value = 7011603
$ws.Range("E200").Value = 45339.375
$ws.Range("F200").Value = "Sabah"
$ws.Range("G200").Value = "FK Kapaz"
$ws.Range("H200").Value = 3
$ws.Range("I200").Value = 2
$ws.Range("J200").Value = "H"
$ws.Range("K200").Value = 1.666
$ws.Range("L200").Value = 3.2
$ws.Range("M200").Value = 5
$ws.Range("N200").Value = 1.666
$ws.Range("O200").Value = 3.2
$ws.Range("P200").Value = 5
$ws.Range("Q200").Value = -0.75
$ws.Range("R200").Value = 1.9
$ws.Range("S200").Value = 1.9
$ws.Range("T200").Value = 2.25
$ws.Range("U200").Value = 1.875
$ws.Range("V200").Value = 1.925
$ws.Range("W200").Value = 0.6659999999999999
$ws.Range("X200").Value = -1
$ws.Range("Y200").Value = -1
$ws.Range("Z200").Value = 0.45
$ws.Range("AA200").Value = -0.5
$ws.Range("AB200").Value = 0.875
$ws.Range("AC200").Value = -1

# ---------------------------------------------------------------------
# 4) Row 201 now carries the final result for match id 7011602
#    (previously this row held the not-yet-played 7011606 fixture)
# ---------------------------------------------------------------------
$ws.Range("B201").Value = 7011602
$ws.Range("E201").Value = 45339.47916666666
$ws.Range("F201").Value = "Sabail FC"
$ws.Range("G201").Value = "Araz FK"
$ws.Range("H201").Value = 2
$ws.Range("I201").Value = 2
$ws.Range("J201").Value = "D"
$ws.Range("K201").Value = 1.8
$ws.Range("L201").Value = 3
$ws.Range("M201").Value = 4.5
$ws.Range("N201").Value = 2.2
$ws.Range("O201").Value = 3.1
$ws.Range("P201").Value = 3
$ws.Range("Q201").Value = -0.25
$ws.Range("R201").Value = 1.975
$ws.Range("S201").Value = 1.825
$ws.Range("T201").Value = 2.25
$ws.Range("U201").Value = 1.9
$ws.Range("V201").Value = 1.9
$ws.Range("W201").Value = -1
$ws.Range("X201").Value = 2.1
$ws.Range("Y201").Value = -1
$ws.Range("Z201").Value = -0.5
$ws.Range("AA201").Value = 0.4125
$ws.Range("AB201").Value = 0.8999999999999999
$ws.Range("AC201").Value = -1

# ---------------------------------------------------------------------
# 5) New row 202: the 7011605 fixture, now played, appended at the end.
#    Copy the A/E number-format styling from row 201 first, then fill in
#    the values.
# ---------------------------------------------------------------------
$ws.Range("A201").Copy()
$ws.Range("A202").PasteSpecial(-4122)
$ws.Range("E201").Copy()
$ws.Range("E202").PasteSpecial(-4122)

$ws.Range("A202").Value = 200
$ws.Range("B202").Value = 7011605
$ws.Range("C202").Value = "Azerbaijan Premier League"
$ws.Range("D202").Value = "Azerbaijan Premier League"
$ws.Range("E202").Value = 45340.35416666666
$ws.Range("F202").Value = "FK Sumqayit"
$ws.Range("G202").Value = "PFK Turan Tovuz"
$ws.Range("H202").Value = 0
$ws.Range("I202").Value = 0
$ws.Range("J202").Value = "D"
$ws.Range("K202").Value = 2.375
$ws.Range("L202").Value = 2.8
$ws.Range("M202").Value = 3
$ws.Range("N202").Value = 2.9
$ws.Range("O202").Value = 2.8
$ws.Range("P202").Value = 2.4
$ws.Range("Q202").Value = 0
$ws.Range("R202").Value = 2
$ws.Range("S202").Value = 1.7
$ws.Range("T202").Value = 2.25
$ws.Range("U202").Value = 1.975
$ws.Range("V202").Value = 1.825
$ws.Range("W202").Value = -1
$ws.Range("X202").Value = 1.8
$ws.Range("Y202").Value = -1
$ws.Range("Z202").Value = 0
$ws.Range("AA202").Value = -0
$ws.Range("AB202").Value = -1
$ws.Range("AC202").Value = 0.825

# ---------------------------------------------------------------------
# 6) New row 203: the 7011606 fixture, now played, appended at the end.
# ---------------------------------------------------------------------
$ws.Range("A201").Copy()
$ws.Range("A203").PasteSpecial(-4122)
$ws.Range("E201").Copy()
$ws.Range("E203").PasteSpecial(-4122)

$ws.Range("A203").Value = 201
$ws.Range("B203").Value = 7011606
$ws.Range("C203").Value = "Azerbaijan Premier League"
$ws.Range("D203").Value = "Azerbaijan Premier League"
$ws.Range("E203").Value = 45340.45833333334
$ws.Range("F203").Value = "Neftchi Baku"
$ws.Range("G203").Value = "Zira IK"
$ws.Range("H203").Value = 1
$ws.Range("I203").Value = 0
$ws.Range("J203").Value = "H"
$ws.Range("K203").Value = 2.2
$ws.Range("L203").Value = 2.9
$ws.Range("M203").Value = 3.2
$ws.Range("N203").Value = 2.1
$ws.Range("O203").Value = 2.75
$ws.Range("P203").Value = 3.75
$ws.Range("Q203").Value = -0.25
$ws.Range("R203").Value = 1.825
$ws.Range("S203").Value = 1.975
$ws.Range("T203").Value = 1.75
$ws.Range("U203").Value = 1.8
$ws.Range("V203").Value = 2
$ws.Range("W203").Value = 1.1
$ws.Range("X203").Value = -1
$ws.Range("Y203").Value = -1
$ws.Range("Z203").Value = 0.825
$ws.Range("AA203").Value = -1
$ws.Range("AB203").Value = -1
$ws.Range("AC203").Value = 1

$wb.Save()
